$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values: safe to assign directly, Excel won't reinterpret them as numbers.
$ws.Range("A1").Value = "Pat"
$ws.Range("C1").Value = "1-027-368-7432"
$ws.Range("A2").Value = "Therese"
$ws.Range("C2").Value = "1-124-506-6318 x44184"

# Numeric-looking text values (e.g. "89067", "1912") must stay text, not become
# real numbers. Writing them straight into B1/B2 via .Value would make Excel
# coerce them to numbers. Instead, write each one into an out-of-the-way
# scratch cell formatted as Text, copy it, and paste-special (values only)
# into the real target cell - that keeps the target cell's original General
# style untouched while still landing a text value. The scratch cell is then
# deleted so no stray content is left behind.

$scratch = $ws.Cells.Item(500, 500)

$scratch.NumberFormat = "@"
$scratch.Value = "89067"
$scratch.Copy()
$ws.Range("B1").PasteSpecial(-4163)  # xlPasteValues
$scratch.Delete()

$scratch = $ws.Cells.Item(500, 500)
$scratch.NumberFormat = "@"
$scratch.Value = "1912"
$scratch.Copy()
$ws.Range("B2").PasteSpecial(-4163)  # xlPasteValues
$scratch.Delete()

$excel.CutCopyMode = 0
